# Update for factor analysis: prepend 22 newer trading-day rows (2018-05-10
# through 2018-06-12) to the M1-volume-share table, pushing the existing
# history down by 22 rows (old row 2 -> new row 24, ... old row 74 -> new row 96).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert 22 blank rows right below the header row (row 1).
$ws.Range("A2:A23").EntireRow.Insert()

# 2) New data for the 22 inserted rows (most-recent date first).
$dates = @("2018-06-12","2018-06-07","2018-06-06","2018-06-05","2018-06-04","2018-06-01","2018-05-31","2018-05-30","2018-05-29","2018-05-28","2018-05-25","2018-05-24","2018-05-23","2018-05-22","2018-05-21","2018-05-18","2018-05-17","2018-05-16","2018-05-15","2018-05-14","2018-05-11","2018-05-10")
$cvals = @(517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99,517035.99)
$dvals = @(3511.10658611,3601.19155035,3599.47311831,3707.38084116,3457.1992198,3755.35287588,4085.28210657,4382.58499717,4076.77096072,3981.53330411,3986.80148632,3822.25584106,4607.57601549,4381.97996424,4612.19518977,3762.12283157,3571.56057785,4018.86335608,3781.14344005,4071.37140022,4141.30873529,4127.18403303)
$evals = @(0.6790835945694998,0.6965069395556004,0.6961745773848353,0.7170450244981978,0.6686573636392313,0.7263233021515582,0.7901349588004503,0.8476363506474666,0.7884888169428979,0.7700688890361385,0.7710878088622033,0.739263013597951,0.891151893602223,0.8475193311475281,0.8920452887177158,0.7276326801873116,0.6907760092000559,0.7772888993820334,0.7313114586181902,0.7874444872241873,0.8009710765569724,0.7982392160031259)

# 3) Make column B hold plain text (so "2018-06-12" style strings aren't
#    auto-converted into date serials) before writing the date strings.
$ws.Range("B2:B23").NumberFormat = "@"

for ($i = 0; $i -lt 22; $i++) {
    $r = 2 + $i
    $ws.Cells.Item($r, 2).Value2 = $dates[$i]
    $ws.Cells.Item($r, 3).Value2 = $cvals[$i]
    $ws.Cells.Item($r, 4).Value2 = $dvals[$i]
    $ws.Cells.Item($r, 5).Value2 = $evals[$i]
}

# 4) Column A: row 2 carries the literal text "0" (matches the existing
#    first-data-row convention), rows 3-23 carry the numeric 0 used by every
#    other data row.
$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value2 = "0"
$ws.Range("A3:A23").Value2 = 0

# 5) The row-insert operation copies the header's bold/border formatting
#    down into the new rows; strip that so the new rows look like ordinary
#    data rows, then restore the same formatting the existing data rows use
#    for column A (bold/border/center, matching the row immediately below
#    the inserted block).
$ws.Range("A2:E23").ClearFormats()

$ws.Range("A24").Copy()
$ws.Range("A2").PasteSpecial(-4122)
$ws.Range("A25").Copy()
$ws.Range("A3:A23").PasteSpecial(-4122)

$excel.CutCopyMode = 0
